$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in new SQL / condition columns (C, D) - written in original authoring order
# so the shared-string table is appended in the same sequence as the source edit.
$ws.Range("C4").Value = 'SELECT count(mt.userPoints) FROM matchtip mt WHERE mt.fk_user = USERID'
$ws.Range("C5").Value = 'SELECT count(mt.userPoints) FROM matchtip mt WHERE mt.fk_user = USERID'
$ws.Range("C6").Value = 'SELECT count(mt.userPoints) FROM matchtip mt WHERE mt.fk_user = USERID'
$ws.Range("C7").Value = 'SELECT count(mt.userPoints) FROM matchtip mt WHERE mt.fk_user = USERID'
$ws.Range("C8").Value = 'SELECT count(mt.userPoints) FROM matchtip mt WHERE mt.fk_user = USERID'
$ws.Range("D4").Value = '>1'
$ws.Range("D5").Value = '>123'
$ws.Range("D6").Value = '>300'
$ws.Range("D7").Value = '>600'
$ws.Range("D8").Value = '>1234'
$ws.Range("D10").Value = 'Matchday mit count = 0'
$ws.Range("D11").Value = 'Matchday mit count = 0'
$ws.Range("C16").Value = 'SELECT md.id, sum(mt.userPoints) FROM matchday md inner join matchtip mt on md.id=mt.fk_match WHERE mt.fk_user = USERID GROUP BY md.id'
$ws.Range("D16").Value = 'Matchday mt summe = 0'
$ws.Range("C12").Value = 'SELECT count(mt.id) from matchtip mt where mt.fk_user=USERID AND mt.userPoints = 4'
$ws.Range("C13").Value = 'SELECT count(mt.id) from matchday md inner join game g on g.matchdayId=md.id inner join matchtip mt on g.id=mt.fk_match where mt.fk_user=USERID AND mt.userPoints = 4 Group by md.id'
$ws.Range("C14").Value = 'SELECT count(mt.id) from matchday md inner join game g on g.matchdayId=md.id inner join matchtip mt on g.id=mt.fk_match where mt.fk_user=USERID AND mt.userPoints = 4 Group by md.id'
$ws.Range("C11").Value = 'SELECT md.id, count(mt.id) FROM  matchday md inner join game g on g.matchdayId=md.id inner join matchtip mt on g.id=mt.fk_match WHERE mt.fk_user = USERID AND mt.userPoints<4 GROUP BY md.id'
$ws.Range("C10").Value = 'SELECT md.id, count(mt.id) FROM  matchday md inner join game g on g.matchdayId=md.id inner join matchtip mt on g.id=mt.fk_match WHERE mt.fk_user = USERID AND mt.userPoints<1 GROUP BY md.id'
$ws.Range("D14").Value = 'count > =5'
$ws.Range("D13").Value = 'count > =3'
$ws.Range("D12").Value = 'count > =1'
$ws.Range("C15").Value = 'Select mt.fk_user, sum(mt.userPoints) from matchtip mt inner join game g on mt.fk_match=g.id inner join matchday md inner join g.matchdayId=md.id group by mt.fk_user where md.id = MATCHDAYID Order by sum(mt.userPoints) Limit 1'
$ws.Range("D15").Value = 'Jeden Spieltag bekommt user mit dieser id das Achievement'
$ws.Range("D20").Value = 'Jeden Spieltag bekommt user mit dieser id das Achievement'
$ws.Range("C2").Value = 'default insert'
$ws.Range("C20").Value = 'Select mt.fk_user, sum(mt.userPoints) from matchtip mt inner join game g on mt.fk_match=g.id inner join matchday md inner join g.matchdayId=md.id inner join league l on md.leagueId=l.id group by mt.fk_user where l.id = leagueId Order by sum(mt.userPoints) Limit 1'
$ws.Range("C18").Value = 'untere abfrage mit allen League ids und dies dann zählen'
$ws.Range("C19").Value = 'untere abfrage mit allen League ids und dies dann zählen'
$ws.Range("C1").Value = 'SQL'
$ws.Range("D1").Value = 'Erfüllt wenn'

# Header cells C1/D1 use the same bold+centered style as A1/B1
$ws.Range("C1:D1").Font.Bold = $true
$ws.Range("C1:D1").HorizontalAlignment = -4108

# Column widths for the two new columns
$ws.Columns.Item(3).ColumnWidth = 243.66666666666666
$ws.Columns.Item(4).ColumnWidth = 54.666666666666664

# View: zoom to 70% and move the selection
$excel.ActiveWindow.Zoom = 70
[void]$ws.Range("B23").Select()
